$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("feature_attrs")

$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 2

$ws.Activate()
$ws.Range("G11").Select()
